$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rectangle 4" holds the big banner text (was "HIP-17" -> "V-Dao Vote")
$bannerShape = $s.Shapes.Item("Rectangle 4")
$bannerRange = $bannerShape.TextFrame.TextRange
[void]$bannerRange.Delete()
[void]$bannerRange.InsertAfter("V-Dao Vote")

# "Rectangle 5" holds the url text (was "talk.harmony.one" -> "gov.harmony.one")
$urlShape = $s.Shapes.Item("Rectangle 5")
$urlShape.TextFrame.TextRange.Text = "gov.harmony.one"
